# Corrected excel sheets for application fix issues
# Applies the corrected loan-schedule figures (after an application fix) to
# the Summary / Repayment schedule / Transactions sheets, and updates the
# saved selection / active-tab state to match.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("NewLoanInput")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsTx = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary.Range("B2").Value = 0
$wsSummary.Range("E2").NumberFormat = "#,##0"
$wsSummary.Range("E2").Value = 10000
$wsSummary.Range("F2").NumberFormat = "#,##0.00"
$wsSummary.Range("F2").Value = 1678.46
# new (blank) cell G2, picking up the column's default formatting
$wsSummary.Range("G2").Value = 0
$wsSummary.Range("G2").ClearContents()

$wsSummary.Range("A3").Value = 578.96
$wsSummary.Range("B3").Value = 0
$wsSummary.Range("E3").Value = 578.96
$wsSummary.Range("F3").Value = 96.98

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$wsSchedule.Range("I2").Value = 5.56
$wsSchedule.Range("K2").Value = 5.56
$wsSchedule.Range("L2").Value = 5.56
# move the trailing blank placeholder cell from P2 to O2
$wsSchedule.Range("N2").Copy() | Out-Null
$wsSchedule.Range("O2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsSchedule.Range("P2").ClearContents()
$wsSchedule.Range("P2").ClearFormats()

# D3 (disbursement date placeholder) and E3 clear back to plain cells
$wsSchedule.Range("A3").Copy() | Out-Null
$wsSchedule.Range("D3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsSchedule.Range("D3").ClearContents()

$wsSchedule.Range("A3").Copy() | Out-Null
$wsSchedule.Range("E3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsSchedule.Range("L3").Value = 0
$wsSchedule.Range("P3").Value = 887.72

$wsSchedule.Range("F5").Value = 864.71
$wsSchedule.Range("G5").Value = 8321.5400000000009
$wsSchedule.Range("H5").Value = 23.01

$wsSchedule.Range("F6").Value = 785.8
$wsSchedule.Range("G6").Value = 7535.74
$wsSchedule.Range("H6").Value = 101.92

$wsSchedule.Range("F7").Value = 813.39
$wsSchedule.Range("G7").Value = 6722.35
$wsSchedule.Range("H7").Value = 74.33

$wsSchedule.Range("F8").Value = 819.21
$wsSchedule.Range("G8").Value = 5903.14
$wsSchedule.Range("H8").Value = 68.510000000000005

$wsSchedule.Range("F9").Value = 829.5
$wsSchedule.Range("G9").Value = 5073.6400000000003
$wsSchedule.Range("H9").Value = 58.22

$wsSchedule.Range("F10").Value = 836.01
$wsSchedule.Range("G10").Value = 4237.63
$wsSchedule.Range("H10").Value = 51.71

$wsSchedule.Range("F11").Value = 844.53
$wsSchedule.Range("G11").Value = 3393.1
$wsSchedule.Range("H11").Value = 43.19

$wsSchedule.Range("F12").Value = 854.25
$wsSchedule.Range("G12").Value = 2538.85
$wsSchedule.Range("H12").Value = 33.47

$wsSchedule.Range("F13").Value = 861.84
$wsSchedule.Range("G13").Value = 1677.01
$wsSchedule.Range("H13").Value = 25.88

$wsSchedule.Range("F14").Value = 871.18
$wsSchedule.Range("G14").Value = 805.83
$wsSchedule.Range("H14").Value = 16.54

$wsSchedule.Range("F15").Value = 805.83
$wsSchedule.Range("H15").Value = 8.2100000000000009
$wsSchedule.Range("K15").Value = 814.04
$wsSchedule.Range("P15").Value = 814.04

# ---------------------------------------------------------------------
# Transactions sheet - figures refreshed after the application fix
# (new transaction ids, recomputed accrual/repayment amounts); the old
# last row (id 15, plain disbursement) is folded into row 7 and the
# physical row 8 is removed.
# ---------------------------------------------------------------------
$wsTx.Range("A2").Value = 3418
$wsTx.Range("D2").Value = "Accrual"
$wsTx.Range("E2").Value = 101.92
$wsTx.Range("G2").Value = 101.92

$wsTx.Range("A3").Value = 3417
$wsTx.Range("D3").Value = "Accrual"
$wsTx.Range("E3").Value = 23.01
$wsTx.Range("G3").Value = 23.01

$wsTx.Range("A4").Value = 3415
$wsTx.Range("J4").NumberFormat = "#,##0"
$wsTx.Range("J4").Value = 10000

$wsTx.Range("A5").Value = 3416
$wsTx.Range("D5").Value = "Accrual"

$wsTx.Range("A6").Value = 3414
$wsTx.Range("C6").Value = 42005
$wsTx.Range("D6").Value = "Repayment (at time of disbursement)"
$wsTx.Range("E6").Value = 5.56
$wsTx.Range("F6").Value = 0
$wsTx.Range("G6").Value = 0
$wsTx.Range("H6").Value = 5.56
$wsTx.Range("J6").NumberFormat = "#,##0"
$wsTx.Range("J6").Value = 5000

$wsTx.Range("A7").Value = 3413
$wsTx.Range("D7").Value = "Disbursement"
$wsTx.Range("E7").NumberFormat = "#,##0"
$wsTx.Range("E7").Value = 5000
$wsTx.Range("H7").Value = 0
$wsTx.Range("K7").Style = "Normal"
$wsTx.Range("L7").Style = "Normal"

# the old row 8 (id 15 / disburse) is now redundant - its data lives in row 7
$wsTx.Rows("8").Delete() | Out-Null

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
$wsTx.Range("B7").Select() | Out-Null
$wsSchedule.Range("E11").Select() | Out-Null
$wsSummary.Range("D4").Select() | Out-Null
